# Auto-generated edit script: applies the commit's cell-value changes
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets of the workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$sheet1Data = @(
  @(21, 8, 4499),
  @(21, 9, 4499),
  @(21, 11, 4499),
  @(21, 13, -4031),
  @(23, 8, 4499),
  @(23, 9, 4499),
  @(23, 11, 4499),
  @(23, 13, -4265),
  @(40, 8, 1929.4),
  @(40, 9, 2375),
  @(40, 10, 1632.3334),
  @(40, 11, 2375),
  @(40, 12, 1632.3334),
  @(40, 13, -2200),
  @(40, 14, -1982.3334),
  @(58, 8, 98.333336),
  @(58, 9, 98.333336),
  @(58, 10, 0),
  @(58, 11, 295.000008),
  @(58, 12, 0),
  @(58, 13, -145.000008),
  @(58, 14, $null),
  @(69, 8, 11000),
  @(69, 10, 11000),
  @(69, 12, 33000),
  @(69, 14, -34748),
  @(72, 8, 11000),
  @(72, 10, 11000),
  @(72, 12, 99000),
  @(72, 14, -107736),
  @(125, 8, $null),
  @(125, 9, $null),
  @(125, 10, $null),
  @(125, 11, $null),
  @(125, 12, $null),
  @(126, 8, $null),
  @(126, 9, $null),
  @(126, 10, $null),
  @(126, 11, $null),
  @(126, 12, $null),
  @(127, 8, $null),
  @(127, 9, $null),
  @(127, 10, $null),
  @(127, 11, $null),
  @(127, 12, $null),
  @(127, 13, $null),
  @(127, 14, $null),
  @(128, 8, $null),
  @(128, 9, $null),
  @(128, 10, $null),
  @(128, 11, $null),
  @(128, 12, $null),
  @(129, 8, $null),
  @(129, 9, $null),
  @(129, 10, $null),
  @(129, 11, $null),
  @(129, 12, $null),
  @(129, 13, $null),
  @(129, 14, $null),
  @(130, 8, $null),
  @(130, 9, $null),
  @(130, 10, $null),
  @(130, 11, $null),
  @(130, 12, $null),
  @(131, 8, $null),
  @(131, 9, $null),
  @(131, 10, $null),
  @(131, 11, $null),
  @(131, 12, $null),
  @(131, 13, $null),
  @(131, 14, $null),
  @(132, 8, $null),
  @(132, 9, $null),
  @(132, 10, $null),
  @(132, 11, $null),
  @(132, 12, $null),
  @(132, 13, $null),
  @(133, 8, $null),
  @(133, 9, $null),
  @(133, 10, $null),
  @(133, 11, $null),
  @(133, 12, $null),
  @(134, 8, $null),
  @(134, 9, $null),
  @(134, 10, $null),
  @(134, 11, $null),
  @(134, 12, $null),
  @(134, 14, $null),
  @(135, 8, $null),
  @(135, 9, $null),
  @(135, 10, $null),
  @(135, 11, $null),
  @(135, 12, $null),
  @(135, 13, $null),
  @(136, 8, $null),
  @(136, 9, $null),
  @(136, 10, $null),
  @(136, 11, $null),
  @(136, 12, $null),
  @(137, 8, $null),
  @(137, 9, $null),
  @(137, 10, $null),
  @(137, 11, $null),
  @(137, 12, $null),
  @(137, 13, $null),
  @(137, 14, $null),
  @(138, 8, $null),
  @(138, 9, $null),
  @(138, 10, $null),
  @(138, 11, $null),
  @(138, 12, $null),
  @(138, 13, $null),
  @(138, 14, $null),
  @(139, 8, $null),
  @(139, 9, $null),
  @(139, 10, $null),
  @(139, 11, $null),
  @(139, 12, $null),
  @(140, 8, $null),
  @(140, 9, $null),
  @(140, 10, $null),
  @(140, 11, $null),
  @(140, 12, $null),
  @(141, 8, $null),
  @(141, 9, $null),
  @(141, 10, $null),
  @(141, 11, $null),
  @(141, 12, $null),
  @(141, 13, $null)
)
foreach ($item in $sheet1Data) {
  $r = $item[0]
  $c = $item[1]
  $v = $item[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets.Item("ARM")
$sheet2Data = @(
  @(122, 8, 0),
  @(122, 9, 0),
  @(122, 11, 0),
  @(122, 13, $null),
  @(132, 8, 8599.667),
  @(132, 9, 4066.1667),
  @(132, 11, 12198.5001),
  @(132, 13, -9668.500100000001)
)
foreach ($item in $sheet2Data) {
  $r = $item[0]
  $c = $item[1]
  $v = $item[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets.Item("BSM")
$sheet3Data = @(
  @(3, 8, 0),
  @(3, 9, 0),
  @(3, 10, 0),
  @(3, 11, 0),
  @(3, 12, 0),
  @(3, 13, $null),
  @(3, 14, $null),
  @(4, 8, 107.14286),
  @(4, 9, 80.333336),
  @(4, 10, 127.25),
  @(4, 11, 80.333336),
  @(4, 12, 127.25),
  @(4, 13, 34.666664),
  @(4, 14, -357.25),
  @(5, 8, 1500),
  @(5, 9, 0),
  @(5, 10, 1500),
  @(5, 11, 0),
  @(5, 12, 1500),
  @(5, 13, $null),
  @(5, 14, -1726),
  @(7, 8, 983.1667),
  @(7, 9, 866.3333),
  @(7, 10, 1100),
  @(7, 11, 866.3333),
  @(7, 12, 1100),
  @(7, 13, -753.3333),
  @(7, 14, -1326),
  @(8, 8, 150),
  @(8, 9, 150),
  @(8, 10, 0),
  @(8, 11, 150),
  @(8, 12, 0),
  @(8, 13, -10),
  @(8, 14, $null),
  @(10, 8, 1631),
  @(10, 9, 1631),
  @(10, 11, 1631),
  @(10, 13, -1491),
  @(11, 8, 389),
  @(11, 9, 100),
  @(11, 10, 461.25),
  @(11, 11, 100),
  @(11, 12, 461.25),
  @(11, 13, 40),
  @(11, 14, -741.25),
  @(12, 8, 675.3333),
  @(12, 9, 1000),
  @(12, 10, 26),
  @(12, 11, 1000),
  @(12, 12, 26),
  @(12, 13, -832),
  @(12, 14, -362),
  @(13, 8, 0),
  @(13, 10, 0),
  @(13, 12, 0),
  @(13, 14, $null),
  @(17, 8, 49),
  @(17, 10, 49),
  @(17, 12, 49),
  @(17, 14, -393),
  @(19, 8, 1281.6),
  @(19, 9, 1009),
  @(19, 10, 1349.75),
  @(19, 11, 1009),
  @(19, 12, 1349.75),
  @(19, 13, -836),
  @(19, 14, -1695.75),
  @(20, 8, 1971.875),
  @(20, 9, 2112.5),
  @(20, 10, 1550),
  @(20, 11, 2112.5),
  @(20, 12, 1550),
  @(20, 13, -1865.5),
  @(20, 14, -2044),
  @(21, 8, 0),
  @(21, 10, 0),
  @(21, 12, 0),
  @(21, 14, $null),
  @(22, 8, 256.25),
  @(22, 9, 274.33334),
  @(22, 10, 202),
  @(22, 11, 274.33334),
  @(22, 12, 202),
  @(22, 13, -101.33334),
  @(22, 14, -548),
  @(23, 8, 7100),
  @(23, 9, 7100),
  @(23, 11, 7100),
  @(23, 13, -6817),
  @(24, 8, 1662.5),
  @(24, 9, 325),
  @(24, 10, 3000),
  @(24, 11, 325),
  @(24, 12, 3000),
  @(24, 13, -90),
  @(24, 14, -3470),
  @(25, 8, 1348.5),
  @(25, 9, 631.3333),
  @(25, 10, 3500),
  @(25, 11, 631.3333),
  @(25, 12, 3500),
  @(25, 13, -396.3333),
  @(25, 14, -3970),
  @(26, 8, 16000),
  @(26, 9, 16000),
  @(26, 10, 0),
  @(26, 11, 16000),
  @(26, 12, 0),
  @(26, 13, -15708),
  @(26, 14, $null),
  @(29, 8, 0),
  @(29, 9, 0),
  @(29, 11, 0),
  @(29, 13, $null),
  @(34, 8, 5000),
  @(34, 10, 5000),
  @(34, 12, 5000),
  @(34, 14, -5228),
  @(35, 8, 33000),
  @(35, 10, 33000),
  @(35, 12, 33000),
  @(35, 14, -33620),
  @(36, 8, 5500),
  @(36, 9, 5500),
  @(36, 11, 5500),
  @(36, 13, -4966),
  @(37, 8, 5481.75),
  @(37, 10, 5481.75),
  @(37, 12, 5481.75),
  @(37, 14, -5755.75),
  @(42, 8, 0),
  @(42, 10, 0),
  @(42, 12, 0),
  @(42, 14, $null),
  @(43, 8, 0),
  @(43, 10, 0),
  @(43, 12, 0),
  @(43, 14, $null),
  @(46, 8, 0),
  @(46, 9, 0),
  @(46, 11, 0),
  @(46, 13, $null),
  @(49, 8, 3610),
  @(49, 10, 3610),
  @(49, 12, 3610),
  @(49, 14, -4088),
  @(50, 8, 59780),
  @(50, 10, 59780),
  @(50, 12, 59780),
  @(50, 14, -60928),
  @(51, 8, 0),
  @(51, 10, 0),
  @(51, 12, 0),
  @(51, 14, $null),
  @(52, 8, 59890.5),
  @(52, 10, 59890.5),
  @(52, 12, 59890.5),
  @(52, 14, -60416.5),
  @(54, 8, 12500),
  @(54, 9, 0),
  @(54, 10, 12500),
  @(54, 11, 0),
  @(54, 12, 12500),
  @(54, 13, $null),
  @(54, 14, -13468),
  @(55, 8, 0),
  @(55, 10, 0),
  @(55, 12, 0),
  @(55, 14, $null),
  @(58, 8, 50999),
  @(58, 10, 50999),
  @(58, 12, 50999),
  @(58, 14, -51587),
  @(60, 8, 100000),
  @(60, 10, 100000),
  @(60, 12, 100000),
  @(60, 14, -101198),
  @(64, 8, 946),
  @(64, 9, 1125.6666),
  @(64, 10, 407),
  @(64, 11, 1125.6666),
  @(64, 12, 407),
  @(64, 13, -900.6666),
  @(64, 14, -857),
  @(67, 8, 946),
  @(67, 9, 1125.6666),
  @(67, 10, 407),
  @(67, 11, 1125.6666),
  @(67, 12, 407),
  @(67, 13, -345.6666),
  @(67, 14, -1967),
  @(75, 8, 10000),
  @(75, 9, 10000),
  @(75, 11, 10000),
  @(75, 13, -9064),
  @(76, 8, 27029.875),
  @(76, 9, 26437.334),
  @(76, 10, 27385.4),
  @(76, 11, 26437.334),
  @(76, 12, 27385.4),
  @(76, 13, -26122.334),
  @(76, 14, -28015.4),
  @(78, 8, 10000),
  @(78, 9, 10000),
  @(78, 11, 30000),
  @(78, 13, -25320),
  @(79, 8, 27029.875),
  @(79, 9, 26437.334),
  @(79, 10, 27385.4),
  @(79, 11, 26437.334),
  @(79, 12, 27385.4),
  @(79, 13, -25345.334),
  @(79, 14, -29569.4),
  @(80, 8, 1330),
  @(80, 9, 745),
  @(80, 10, 2500),
  @(80, 11, 745),
  @(80, 12, 2500),
  @(80, 13, 253),
  @(80, 14, -4496),
  @(81, 8, 44561.5),
  @(81, 10, 44561.5),
  @(81, 12, 44561.5),
  @(81, 14, -46683.5),
  @(82, 8, 5750),
  @(82, 9, 5750),
  @(82, 10, 0),
  @(82, 11, 5750),
  @(82, 12, 0),
  @(82, 13, -5367),
  @(82, 14, $null),
  @(83, 8, 1330),
  @(83, 9, 745),
  @(83, 10, 2500),
  @(83, 11, 3725),
  @(83, 12, 12500),
  @(83, 13, 1267),
  @(83, 14, -22484),
  @(84, 8, 44561.5),
  @(84, 10, 44561.5),
  @(84, 12, 133684.5),
  @(84, 14, -144292.5),
  @(85, 8, 5750),
  @(85, 9, 5750),
  @(85, 10, 0),
  @(85, 11, 5750),
  @(85, 12, 0),
  @(85, 13, -4424),
  @(85, 14, $null),
  @(86, 8, 2028.5714),
  @(86, 9, 2028.5714),
  @(86, 11, 2028.5714),
  @(86, 13, -905.5714),
  @(88, 8, 15799.8),
  @(88, 10, 15799.8),
  @(88, 12, 15799.8),
  @(88, 14, -16611.8),
  @(89, 8, 2028.5714),
  @(89, 9, 2028.5714),
  @(89, 11, 10142.857),
  @(89, 13, -4526.857),
  @(91, 8, 15799.8),
  @(91, 10, 15799.8),
  @(91, 12, 15799.8),
  @(91, 14, -18607.8),
  @(92, 8, 19664.334),
  @(92, 10, 19664.334),
  @(92, 12, 19664.334),
  @(92, 14, -24656.334),
  @(94, 8, 747.8889),
  @(94, 9, 747.8889),
  @(94, 11, 747.8889),
  @(94, 13, -296.8889),
  @(95, 8, 14399.5),
  @(95, 10, 14399.5),
  @(95, 12, 14399.5),
  @(95, 14, -19891.5),
  @(96, 8, 20000),
  @(96, 9, 20000),
  @(96, 11, 20000),
  @(96, 13, -17254),
  @(97, 8, 0),
  @(97, 9, 0),
  @(97, 11, 0),
  @(97, 13, $null),
  @(98, 8, 70542),
  @(98, 10, 70542),
  @(98, 12, 70542),
  @(98, 14, -76532),
  @(99, 8, 1399.7778),
  @(99, 9, 1399.7778),
  @(99, 11, 1399.7778),
  @(99, 13, 98.22219999999993),
  @(100, 8, 22300),
  @(100, 10, 22300),
  @(100, 12, 22300),
  @(100, 14, -24464),
  @(102, 8, 47500),
  @(102, 9, 20000),
  @(102, 10, 75000),
  @(102, 11, 20000),
  @(102, 12, 75000),
  @(102, 13, -16755),
  @(102, 14, -81490),
  @(103, 8, 20007.2),
  @(103, 10, 20007.2),
  @(103, 12, 20007.2),
  @(103, 14, -22351.2),
  @(105, 8, 2186.2778),
  @(105, 9, 1873.1111),
  @(105, 10, 2499.4443),
  @(105, 11, 1873.1111),
  @(105, 12, 2499.4443),
  @(105, 13, -126.1111000000001),
  @(105, 14, -5993.4443),
  @(106, 8, 0),
  @(106, 10, 0),
  @(106, 12, 0),
  @(106, 14, $null),
  @(107, 8, 1544.3334),
  @(107, 9, 1544.3334),
  @(107, 10, 0),
  @(107, 11, 1544.3334),
  @(107, 12, 0),
  @(107, 13, 375.6666),
  @(107, 14, $null),
  @(109, 8, 0),
  @(109, 10, 0),
  @(109, 12, 0),
  @(109, 14, $null),
  @(120, 8, 49999.5),
  @(120, 10, 49999.5),
  @(120, 12, 49999.5),
  @(120, 14, -59675.5),
  @(121, 8, 59890.5),
  @(121, 10, 59890.5),
  @(121, 12, 59890.5),
  @(121, 14, -63384.5),
  @(122, 8, 0),
  @(122, 10, 0),
  @(122, 12, 0),
  @(122, 14, $null),
  @(130, 8, 0),
  @(130, 10, 0),
  @(130, 12, 0),
  @(130, 14, $null),
  @(134, 8, 9187.5),
  @(134, 9, 8000),
  @(134, 10, 9583.333),
  @(134, 11, 24000),
  @(134, 12, 28749.999),
  @(134, 13, -21465),
  @(134, 14, -33819.999),
  @(135, 8, 80000),
  @(135, 10, 80000),
  @(135, 12, 80000),
  @(135, 14, -90140),
  @(139, 8, 0),
  @(139, 9, 0),
  @(139, 11, 0),
  @(139, 13, $null),
  @(140, 8, 0),
  @(140, 10, 0),
  @(140, 12, 0),
  @(140, 14, $null)
)
foreach ($item in $sheet3Data) {
  $r = $item[0]
  $c = $item[1]
  $v = $item[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets.Item("CRP")
$sheet4Data = @(
  @(31, 8, 4623.9473),
  @(31, 9, 3542.2727),
  @(31, 10, 6111.25),
  @(31, 11, 3542.2727),
  @(31, 12, 6111.25),
  @(31, 13, -3247.2727),
  @(31, 14, -6701.25),
  @(34, 8, 4623.9473),
  @(34, 9, 3542.2727),
  @(34, 10, 6111.25),
  @(34, 11, 3542.2727),
  @(34, 12, 6111.25),
  @(34, 13, -3340.2727),
  @(34, 14, -6515.25)
)
foreach ($item in $sheet4Data) {
  $r = $item[0]
  $c = $item[1]
  $v = $item[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets.Item("CUL")
$sheet5Data = @(
  @(11, 8, 652),
  @(11, 10, 716.6667),
  @(11, 12, 2150.0001),
  @(11, 14, -2430.0001)
)
foreach ($item in $sheet5Data) {
  $r = $item[0]
  $c = $item[1]
  $v = $item[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets.Item("GSM")
$sheet6Data = @(
  @(25, 8, 4989),
  @(25, 10, 4989),
  @(25, 12, 4989),
  @(25, 14, -6047),
  @(97, 8, 699.6667),
  @(97, 9, 0),
  @(97, 11, 0),
  @(97, 12, 699.6667),
  @(97, 13, $null),
  @(97, 14, -1691.6667),
  @(122, 8, 998.6667),
  @(122, 9, 998.6667),
  @(122, 11, 2996.0001),
  @(122, 13, -546.0001000000002),
  @(124, 8, 0),
  @(124, 10, 0),
  @(124, 12, 0),
  @(124, 14, $null),
  @(125, 8, $null),
  @(125, 9, $null),
  @(125, 10, $null),
  @(125, 11, $null),
  @(125, 12, $null),
  @(126, 8, $null),
  @(126, 9, $null),
  @(126, 10, $null),
  @(126, 11, $null),
  @(126, 12, $null),
  @(127, 8, $null),
  @(127, 9, $null),
  @(127, 10, $null),
  @(127, 11, $null),
  @(127, 12, $null),
  @(128, 8, $null),
  @(128, 9, $null),
  @(128, 10, $null),
  @(128, 11, $null),
  @(128, 12, $null),
  @(129, 8, $null),
  @(129, 9, $null),
  @(129, 10, $null),
  @(129, 11, $null),
  @(129, 12, $null),
  @(130, 8, $null),
  @(130, 9, $null),
  @(130, 10, $null),
  @(130, 11, $null),
  @(130, 12, $null),
  @(130, 14, $null),
  @(131, 8, $null),
  @(131, 9, $null),
  @(131, 10, $null),
  @(131, 11, $null),
  @(131, 12, $null),
  @(131, 14, $null),
  @(132, 8, $null),
  @(132, 9, $null),
  @(132, 10, $null),
  @(132, 11, $null),
  @(132, 12, $null),
  @(132, 13, $null),
  @(132, 14, $null),
  @(133, 8, $null),
  @(133, 9, $null),
  @(133, 10, $null),
  @(133, 11, $null),
  @(133, 12, $null),
  @(133, 14, $null),
  @(134, 8, $null),
  @(134, 9, $null),
  @(134, 10, $null),
  @(134, 11, $null),
  @(134, 12, $null),
  @(134, 14, $null),
  @(135, 8, $null),
  @(135, 9, $null),
  @(135, 10, $null),
  @(135, 11, $null),
  @(135, 12, $null),
  @(136, 8, $null),
  @(136, 9, $null),
  @(136, 10, $null),
  @(136, 11, $null),
  @(136, 12, $null),
  @(136, 14, $null),
  @(137, 8, $null),
  @(137, 9, $null),
  @(137, 10, $null),
  @(137, 11, $null),
  @(137, 12, $null),
  @(138, 8, $null),
  @(138, 9, $null),
  @(138, 10, $null),
  @(138, 11, $null),
  @(138, 12, $null),
  @(139, 8, $null),
  @(139, 9, $null),
  @(139, 10, $null),
  @(139, 11, $null),
  @(139, 12, $null),
  @(140, 8, $null),
  @(140, 9, $null),
  @(140, 10, $null),
  @(140, 11, $null),
  @(140, 12, $null),
  @(140, 14, $null),
  @(141, 8, $null),
  @(141, 9, $null),
  @(141, 10, $null),
  @(141, 11, $null),
  @(141, 12, $null),
  @(141, 14, $null)
)
foreach ($item in $sheet6Data) {
  $r = $item[0]
  $c = $item[1]
  $v = $item[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

$ws = $wb.Worksheets.Item("LTW")
$sheet7Data = @(
  @(22, 8, 894.1429),
  @(22, 10, 915),
  @(22, 12, 915),
  @(22, 14, -1505),
  @(27, 8, 894.1429),
  @(27, 10, 915),
  @(27, 12, 915),
  @(27, 14, -1129),
  @(46, 8, 4599.8335),
  @(46, 9, 4424.5),
  @(46, 10, 4687.5),
  @(46, 11, 4424.5),
  @(46, 12, 4687.5),
  @(46, 13, -4236.5),
  @(46, 14, -5063.5),
  @(68, 8, 6713.7144),
  @(68, 9, 6749),
  @(68, 10, 6666.6665),
  @(68, 11, 6749),
  @(68, 12, 6666.6665),
  @(68, 13, -6000),
  @(68, 14, -8164.6665),
  @(71, 8, 6713.7144),
  @(71, 9, 6749),
  @(71, 10, 6666.6665),
  @(71, 11, 33745),
  @(71, 12, 33333.3325),
  @(71, 13, -30001),
  @(71, 14, -40821.3325),
  @(100, 8, 1866.3334),
  @(100, 9, 1866.3334),
  @(100, 11, 1866.3334),
  @(100, 13, -1325.3334),
  @(101, 8, 13989.833),
  @(101, 10, 13989.833),
  @(101, 12, 13989.833),
  @(101, 14, -20479.833),
  @(124, 8, $null),
  @(124, 9, $null),
  @(124, 10, $null),
  @(124, 11, $null),
  @(124, 12, $null),
  @(125, 8, $null),
  @(125, 9, $null),
  @(125, 10, $null),
  @(125, 11, $null),
  @(125, 12, $null),
  @(127, 8, $null),
  @(127, 9, $null),
  @(127, 10, $null),
  @(127, 11, $null),
  @(127, 12, $null),
  @(128, 8, $null),
  @(128, 9, $null),
  @(128, 10, $null),
  @(128, 11, $null),
  @(128, 12, $null),
  @(129, 8, $null),
  @(129, 9, $null),
  @(129, 10, $null),
  @(129, 11, $null),
  @(129, 12, $null),
  @(130, 8, $null),
  @(130, 9, $null),
  @(130, 10, $null),
  @(130, 11, $null),
  @(130, 12, $null),
  @(130, 14, $null),
  @(131, 8, $null),
  @(131, 9, $null),
  @(131, 10, $null),
  @(131, 11, $null),
  @(131, 12, $null),
  @(131, 13, $null),
  @(131, 14, $null),
  @(132, 8, $null),
  @(132, 9, $null),
  @(132, 10, $null),
  @(132, 11, $null),
  @(132, 12, $null),
  @(132, 13, $null),
  @(132, 14, $null),
  @(133, 8, $null),
  @(133, 9, $null),
  @(133, 10, $null),
  @(133, 11, $null),
  @(133, 12, $null),
  @(133, 14, $null),
  @(134, 8, $null),
  @(134, 9, $null),
  @(134, 10, $null),
  @(134, 11, $null),
  @(134, 12, $null),
  @(134, 13, $null),
  @(135, 8, $null),
  @(135, 9, $null),
  @(135, 10, $null),
  @(135, 11, $null),
  @(135, 12, $null),
  @(136, 8, $null),
  @(136, 9, $null),
  @(136, 10, $null),
  @(136, 11, $null),
  @(136, 12, $null),
  @(136, 13, $null),
  @(136, 14, $null),
  @(137, 8, $null),
  @(137, 9, $null),
  @(137, 10, $null),
  @(137, 11, $null),
  @(137, 12, $null),
  @(138, 8, $null),
  @(138, 9, $null),
  @(138, 10, $null),
  @(138, 11, $null),
  @(138, 12, $null),
  @(138, 14, $null),
  @(139, 8, $null),
  @(139, 9, $null),
  @(139, 10, $null),
  @(139, 11, $null),
  @(139, 12, $null),
  @(139, 13, $null),
  @(139, 14, $null),
  @(140, 8, $null),
  @(140, 9, $null),
  @(140, 10, $null),
  @(140, 11, $null),
  @(140, 12, $null),
  @(140, 13, $null),
  @(140, 14, $null),
  @(141, 8, $null),
  @(141, 9, $null),
  @(141, 10, $null),
  @(141, 11, $null),
  @(141, 12, $null)
)
foreach ($item in $sheet7Data) {
  $r = $item[0]
  $c = $item[1]
  $v = $item[2]
  if ($null -eq $v) {
    $ws.Cells.Item($r, $c).Value = ""
  } else {
    $ws.Cells.Item($r, $c).Value = $v
  }
}

